$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.527.73'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.62%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.577.08'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.29%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.13%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '586.86'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.02%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.07'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.18%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.587'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.68%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.107'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.89%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.60'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.93%  '

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.25%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.351'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.79%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.26'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.36%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.031.38'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.06%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '63.337.78'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.48%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000146'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.09%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.569.10'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.39%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.10'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -2.37%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '341.72'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.54%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.32'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.99%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.62'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -3.66%  '

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.16%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.81'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.80%  '

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +4.93%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.62'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.33%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.164'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -3.72%  '

$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'Aptos'
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.97'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -3.35%  '

$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.19%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.23'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -2.59%  '

$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.23%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '474.05'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +2.17%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0₃0799'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.44%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.68'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +3.52%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '176.22'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.25%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.06%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.397'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -2.33%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.86'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -1.85%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.55'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.65%  '

$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.74'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.18%  '

$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'USDe'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.997'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.23%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '40.08'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.41%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '158.27'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +4.62%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.71'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -3.15%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '21.66'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +2.80%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.633'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +3.08%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0537'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.51%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0962'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.37%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0236'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.98%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '17.98'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.56%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '11.37'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.30%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.69'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -3.44%  '
